$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.01
$ws.Range("C2").Value = 2.93
$ws.Range("D2").Value = 0.83
$ws.Range("E2").Value = 1.85
$ws.Range("G2").Value = 1.08

# Row 3
$ws.Range("C3").Value = 2.16
$ws.Range("D3").Value = 0.83
$ws.Range("E3").Value = 1.85
$ws.Range("G3").Value = 1.64

# Row 4
$ws.Range("D4").Value = 0.83
$ws.Range("E4").Value = 1.85
